$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 396: revised input values (cumulative column B recalculates via its shared formula) ---
$ws.Range("C396").Value = 43

# --- Row 398: revised input values ---
$ws.Range("C398").Value = 124

# --- Row 399: revised input values ---
$ws.Range("C399").Value = 70

# --- Row 400: newly-reported day, fill in the data (previously blank placeholder row) ---
$ws.Range("C400").Value = 16
$ws.Range("E400").Value = 6
$ws.Range("F400").Value = 6
$ws.Range("G400").Value = 22

# L400 / M400 are formatted with a Text ("@") number format, so a plain
# Range.Value assignment there would be stored as a text string (matching
# real Excel's "typed into a text cell" behaviour). The source data for
# this row is numeric like all the other rows, so borrow the General
# formatting from a plain numeric cell long enough to write true numbers,
# then restore each cell's own original (Text-formatted) style by copying
# it back from an identically-styled neighbour - this keeps the style
# table untouched while the stored values become real numbers.
$ws.Range("C400").Copy()
$ws.Range("L400").PasteSpecial(-4122)
$ws.Range("L400").Value = 0
$ws.Range("L396").Copy()
$ws.Range("L400").PasteSpecial(-4122)

$ws.Range("C400").Copy()
$ws.Range("M400").PasteSpecial(-4122)
$ws.Range("M400").Value = 0
$ws.Range("M396").Copy()
$ws.Range("M400").PasteSpecial(-4122)
